# Apply targeted odds updates to Sheet1 as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "M6"  = 1.06
    "N6"  = 10
    "Q8"  = 1.88
    "R8"  = 1.98
    "G9"  = 1.45
    "M9"  = 1.07
    "N9"  = 8.5
    "W9"  = 6
    "G10" = 1.7
    "I10" = 5.5
    "J10" = 2.38
    "G12" = 2.72
    "I12" = 3
    "J12" = 3.35
    "K12" = 1.83
    "L12" = 3.6
    "O12" = 1.47
    "P12" = 2.32
    "Q12" = 2.4
    "S12" = 1.52
    "T12" = 2.2
    "U12" = 1.9
    "V12" = 1.72
    "AC12" = 5.9
    "AH12" = 7
    "AI12" = 14.5
    "AJ12" = 10.75
    "AT12" = 2.18
    "BA12" = 120
    "BB12" = 350
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
